# Insert a new data row before the current row 111 ("Start Ruby" / "Primera",
# 2021-08-05 record), shifting every subsequent row down by one. The new row
# carries a fresh weekly price observation for 2022-01-21.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(111).Insert()

$ws.Range("A111").Value = 10
$ws.Range("B111").Value = "Vega Modelo de Temuco"
$ws.Range("C111").Value = "La Araucanía"
$ws.Range("D111").Value = "2022-01-21"
$ws.Range("E111").Value = 9
$ws.Range("F111").Value = "Fruta"
$ws.Range("G111").Value = 100102
$ws.Range("H111").Value = "Cítricos"
$ws.Range("I111").Value = 100102006
$ws.Range("J111").Value = "Pomelo"
$ws.Range("K111").Value = "Start Ruby"
$ws.Range("L111").Value = "Primera"
$ws.Range("M111").Value = 55
$ws.Range("N111").Value = 15000
$ws.Range("O111").Value = 15000
$ws.Range("P111").Value = 15000
$ws.Range("Q111").Value = "`$/bandeja 15 kilos empedrada"
$ws.Range("R111").Value = "Región de O'Higgins"
$ws.Range("S111").Value = 1000
$ws.Range("T111").Value = 15
